$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from the end of the "...und den Zusammen"
#    paragraph to the (now) empty "Listenabsatz" paragraph right after the
#    "Im Unterschenkel..." paragraph.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$targetPara = $d.Paragraphs(11)
$d.Bookmarks.Add("_GoBack", $targetPara.Range)

# 2) Trim the trailing space in the "Ja -> MS ..." list item.
$d.Content.Find.Execute(
    "Ja -> MS – Studie muss durch Segmentierungsalgorithmus segmentiert werden (fertig?) ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ja -> MS – Studie muss durch Segmentierungsalgorithmus segmentiert werden (fertig?)",
    2) | Out-Null

# 3) Add a new list item right after it with the "Radiomics ..." note.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Radiomics</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> hat in der </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Structure</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> 42 Werte, obwohl es 57 Features zur Auswahl gibt</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>in der Pipeline. Das gleiche Problem, mit anderen zahlen ist bei PORTS</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>-&gt; Beantwortet!!!!</w:t></w:r>' +
    '</w:p>'

$insertionPoint.InsertXML($newParaXml)
